$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, shifting existing rows 75-142 down to 76-143
$ws.Rows("75:75").Insert()

# Populate the newly inserted row 75 with the new weekly record
$ws.Cells.Item(75, 1).Value = 6
$ws.Cells.Item(75, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(75, 3).Value = "Metropolitana"
$ws.Cells.Item(75, 4).Value = 44546
$ws.Cells.Item(75, 5).Value = 13
$ws.Cells.Item(75, 6).Value = 100112001
$ws.Cells.Item(75, 7).Value = "Berenjena"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 180
$ws.Cells.Item(75, 11).Value = 7000
$ws.Cells.Item(75, 12).Value = 8000
$ws.Cells.Item(75, 13).Value = 7556
$ws.Cells.Item(75, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(75, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(75, 16).Value = 151
$ws.Cells.Item(75, 17).Value = 50
$ws.Cells.Item(75, 18).Value = "Hortaliza"
